$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# "Volume 31   Number  10" -> "Volume 31   Number  11"
$ws.Range("A8").Value = "Volume 31   Number  11"
# "Report Covering the Week  3/4/2024  Through  3/10/2024" -> "...3/11/2024...3/17/2024"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# --- Crime statistics grid updates (Week to Date / 28 Day / YTD / % changes) ---
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("N14").Value = -80
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("L15").Value = -60
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = 45.714285714285
$ws.Range("L16").Value = -15
$ws.Range("M16").Value = 6.25
$ws.Range("N16").Value = -72.872340425531
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 35.714285714285
$ws.Range("I17").Value = 102
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = 8.510638297872
$ws.Range("L17").Value = 15.909090909090
$ws.Range("M17").Value = 72.881355932203
$ws.Range("N17").Value = -35.849056603773
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 125
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 18.181818181818
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 44
$ws.Range("K18").Value = -29.545454545454
$ws.Range("L18").Value = -40.384615384615
$ws.Range("M18").Value = 29.166666666666
$ws.Range("N18").Value = -84.878048780487
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -9.677419354838
$ws.Range("I19").Value = 61
$ws.Range("J19").Value = 62
$ws.Range("K19").Value = -1.612903225806
$ws.Range("L19").Value = -16.438356164383
$ws.Range("M19").Value = 12.962962962963
$ws.Range("N19").Value = -28.235294117647
$ws.Range("C20").Value = 1
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 20
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = -4.761904761904
$ws.Range("M20").Value = 122.222222222222
$ws.Range("N20").Value = -70.149253731343
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 105.263157894737
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = 6.25
$ws.Range("I21").Value = 269
$ws.Range("J21").Value = 255
$ws.Range("K21").Value = 5.490196078431
$ws.Range("L21").Value = -11.221122112211
$ws.Range("M21").Value = 33.830845771144
$ws.Range("N21").Value = -63.100137174211
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("M22").Value = 33.333333333333
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 700
$ws.Range("F23").Value = 22
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 46.666666666666
$ws.Range("I23").Value = 55
$ws.Range("J23").Value = 44
$ws.Range("K23").Value = 25
$ws.Range("L23").Value = 12.244897959183
$ws.Range("M23").Value = 66.666666666666
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -15.384615384615
$ws.Range("F24").Value = 48
$ws.Range("H24").Value = -15.789473684210
$ws.Range("I24").Value = 160
$ws.Range("J24").Value = 175
$ws.Range("K24").Value = -8.571428571428
$ws.Range("L24").Value = -8.571428571428
$ws.Range("M24").Value = 23.076923076923
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = -25
$ws.Range("I25").Value = 36
$ws.Range("J25").Value = 38
$ws.Range("K25").Value = -5.263157894736
$ws.Range("L25").Value = -14.285714285714
$ws.Range("C26").Value = 17
$ws.Range("E26").Value = 70
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 25.641025641025
$ws.Range("I26").Value = 118
$ws.Range("J26").Value = 117
$ws.Range("K26").Value = 0.854700854700
$ws.Range("L26").Value = 13.461538461538
$ws.Range("M26").Value = -31.395348837209
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("L27").Value = -42.857142857142
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 62.5
$ws.Range("L28").Value = -7.142857142857
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -50
$ws.Range("M29").Value = 40
$ws.Range("N29").Value = -69.565217391304
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -50
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = -80
